# [Fonds de solidarite] Add 2021-01-22 data
#
# The sheet stores every value (even the numeric-looking ones) as literal
# text (inline strings with no leading apostrophe). Plain Excel "Value ="
# assignment of a numeric-looking string auto-converts it to a real number,
# which would change the cell's stored type. Prefixing with a single quote
# forces Excel to keep/store it as text, matching the source data's type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auvergne-Rhône-Alpes / SARL (row 4)
$ws.Range("C4").Value = "'1576"
$ws.Range("D4").Value = "'1374"
$ws.Range("E4").Value = "'13241924.70"

# Auvergne-Rhône-Alpes / SAS (row 6)
$ws.Range("C6").Value = "'1092"
$ws.Range("D6").Value = "'968"
$ws.Range("E6").Value = "'7840141.57"

# Centre-Val de Loire / SARL (row 23)
$ws.Range("C23").Value = "'481"
$ws.Range("E23").Value = "'3962452.17"

# Centre-Val de Loire / SAS (row 24)
$ws.Range("C24").Value = "'215"
$ws.Range("E24").Value = "'1497756.95"

# Grand Est / SAS (row 36)
$ws.Range("C36").Value = "'573"
$ws.Range("E36").Value = "'3858033.56"

# Île-de-France / SARL (row 58)
$ws.Range("C58").Value = "'6904"
$ws.Range("D58").Value = "'6511"
$ws.Range("E58").Value = "'37769652.78"

# Île-de-France / SA a conseil d'administration (row 59)
$ws.Range("C59").Value = "'27"
$ws.Range("E59").Value = "'397400.20"

# Île-de-France / SAS (row 60)
$ws.Range("C60").Value = "'6751"
$ws.Range("E60").Value = "'31094761.20"

# Nouvelle-Aquitaine / Entrepreneur individuel (row 79)
$ws.Range("C79").Value = "'456"
$ws.Range("E79").Value = "'1490029.96"
